# This edit cyclically rotates the data held in rows 9-13 of the active
# sheet. Content that is currently in row 9 ends up in row 10, row 10's
# content ends up in row 13, row 13's content ends up in row 11, row 11's
# content ends up in row 12, and row 12's content ends up back in row 9
# (a single 5-cycle: 9 -> 10 -> 13 -> 11 -> 12 -> 9). Nothing else on the
# sheet is touched.
#
# Only the column groups that are actually populated for these rows are
# moved (K:O, X, AF, AH, AJ:AS, AU:AV are always blank for rows 9-13, so
# they are left alone; AT and AY are always-empty placeholder cells in
# every one of these rows both before and after the edit, so they are left
# untouched too). Columns I, Y, Z, AA, AB hold text that looks like a
# number / date ("1", "2023-08-11", "00:00", ...) - a plain Value2 write
# lets Excel's COM layer silently reinterpret such text as a real
# number/date, so those destination cells are pre-formatted as Text ("@")
# before the write and reset back to the Normal style afterwards so no
# stray number format sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column-letter groups (start,end) that hold data for rows 9:13, keyed by a
# short name so we can stash/retrieve them from a hashtable per row.
$colGroups = @{
    "AJ" = @("A", "J")
    "PW" = @("P", "W")
    "YAE" = @("Y", "AE")
    "AG" = @("AG", "AG")
    "AI" = @("AI", "AI")
    "AWAX" = @("AW", "AX")
}

function Get-RowValues($row) {
    $vals = @{}
    foreach ($key in $colGroups.Keys) {
        $grp = $colGroups[$key]
        $rng = "$($grp[0])$row" + ":" + "$($grp[1])$row"
        $vals[$key] = $ws.Range($rng).Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    # Force text-prone columns to Text format first so date-/number-looking
    # strings ("1", "2023-08-11", "00:00") round-trip as text, matching the
    # source workbook (which stores them as inline strings).
    $textRng1 = "I$row" + ":I$row"
    $textRng2 = "Y$row" + ":AB$row"
    $ws.Range($textRng1).NumberFormat = "@"
    $ws.Range($textRng2).NumberFormat = "@"

    foreach ($key in $colGroups.Keys) {
        $grp = $colGroups[$key]
        $rng = "$($grp[0])$row" + ":" + "$($grp[1])$row"
        $ws.Range($rng).Value2 = $vals[$key]
    }

    # Drop the temporary Text number format again so the cells end up
    # style-less, like the originals.
    $ws.Range($textRng1).Style = "Normal"
    $ws.Range($textRng2).Style = "Normal"
}

# Snapshot every source row before any writes happen.
$row9  = Get-RowValues 9
$row10 = Get-RowValues 10
$row11 = Get-RowValues 11
$row12 = Get-RowValues 12
$row13 = Get-RowValues 13

# Apply the cycle: 9<-12, 12<-11, 11<-13, 13<-10, 10<-(old 9)
Set-RowValues 9  $row12
Set-RowValues 12 $row11
Set-RowValues 11 $row13
Set-RowValues 13 $row10
Set-RowValues 10 $row9
